$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (H1:M1)
$ws.Range("H1").Value = "ubicacion_deuda"
$ws.Range("I1").Value = "nombre_deuda"
$ws.Range("J1").Value = "ubicacion_vencimientos"
$ws.Range("K1").Value = "nombre_vencimientos"
$ws.Range("L1").Value = "ubicacion_ddjj"
$ws.Range("M1").Value = "nombre_ddjj"

# Mirror the existing header style (bold font, border, centered/top alignment) from G1
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:M1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Row 2 data
$ws.Range("H2").Value = "./Descargas"
$ws.Range("I2").Value = "deuda-demo"
$ws.Range("J2").Value = "./Descargas"
$ws.Range("K2").Value = "vencimientos-demo"
$ws.Range("L2").Value = "./Descargas"
$ws.Range("M2").Value = "ddjj-demo"

# Row 3 data
$ws.Range("H3").Value = "./Descargas"
$ws.Range("I3").Value = "deuda-no"
$ws.Range("J3").Value = "./Descargas"
$ws.Range("K3").Value = "vencimientos-no"
$ws.Range("L3").Value = "./Descargas"
$ws.Range("M3").Value = "ddjj-no"
